$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos text replaced with docente name ---
$ws.Range("B10").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C10").Value = '519033 - Carlos Yujiro Shigue'

# --- Row 13: becomes 'Programa resumido:' row, B/C copy the date text from B8/C8 ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4163) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4163) | Out-Null

# --- Row 14: becomes 'Short syllabus:' (A only, B/C cleared) ---
$ws.Rows.Item(14).RowHeight = 60
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# --- Row 15: becomes 'Programa:' with docente name in B/C ---
$ws.Rows.Item(15).RowHeight = 120
$ws.Range("A15").Value = 'Programa:'
$ws.Range("B15").Value = '519033 - Carlos Yujiro Shigue'
$ws.Range("C15").Value = '519033 - Carlos Yujiro Shigue'

# --- Row 16: becomes 'Syllabus:' ---
$ws.Rows.Item(16).RowHeight = 120
$ws.Range("A16").Value = 'Syllabus:'

# --- Row 17: becomes 'Avaliação:' (A only, B/C cleared, default height) ---
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows.Item(17).AutoFit()

# --- Row 18: becomes 'Método:' with B/C filled from the docente name row ---
# B18/C18 are brand-new cells; paste the number formats from an existing B/C
# pair first so they pick up style index 2/3 instead of the default column style.
$ws.Rows.Item(18).RowHeight = 60
$ws.Range("A18").Value = 'Método:'
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B18").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C18").PasteSpecial(-4122) | Out-Null
$ws.Range("B18").Value = '5817692 - Katia Cristiane Gandolpho Candioto'
$ws.Range("C18").Value = '5817692 - Katia Cristiane Gandolpho Candioto'

# --- Row 19: becomes 'Critério:' with B/C filled ---
$ws.Rows.Item(19).RowHeight = 60
$ws.Range("A19").Value = 'Critério:'
$ws.Range("B9").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C19").PasteSpecial(-4122) | Out-Null
$ws.Range("B19").Value = 'Aulas expositivas, trabalhos e aulas práticas. Aulas com softwares para desenho técnico.'
$ws.Range("C19").Value = 'Aulas expositivas, trabalhos e aulas práticas. Aulas com softwares para desenho técnico.'

# --- Row 20: becomes 'Norma de recuperação:' ---
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("B20").Value = 'Média aritmética das notas de atividades em aula e extra aula.'
$ws.Range("C20").Value = 'Média aritmética das notas de atividades em aula e extra aula.'

# --- Row 21: becomes 'Bibliografia:' ---
$ws.Rows.Item(21).RowHeight = 120
$ws.Range("A21").Value = 'Bibliografia:'
$ws.Range("B21").Value = 'Devido às características práticas da disciplina, não será oferecida recuperação'
$ws.Range("C21").Value = 'Devido às características práticas da disciplina, não será oferecida recuperação'

# --- Remove now-obsolete trailing rows 22 and 23 ---
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()

